$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data (rows 2-11). Column B holds date-text strings and column E
# holds zero-padded codes that must remain plain text (not be silently
# reinterpreted as a date serial / number by Excel).
$data = @(
    @{ Row=2;  A=1; B="2025-04-28"; C=250; D="MAP SERVICOS DE CONSERVACAO - EIRELI";               E="000098"; F="PANO DE CHAO FLANELADO C REFORCADO ITATEX 42x62CM";       G=-15;  H=$false },
    @{ Row=3;  A=2; B="2025-04-28"; C=60;  D="MAP SERVICOS DE CONSERVACAO - EIRELI";               E="000347"; F="INSETICIDA BUZZOFF AEROSSOL 300ML";                        G=2;    H=$true  },
    @{ Row=4;  A=5; B="2025-04-28"; C=70;  D="MAP SERVICOS DE CONSERVACAO - EIRELI";               E="000349"; F="DESODORISADOR LADY AEROSSOL 360ML TALCO SUAVE CARINHO";   G=344;  H=$true  },
    @{ Row=5;  A=7; B="2025-04-28"; C=20;  D="LUCAS CLIENTE NOVO";                                 E="000158"; F="AZULIM LIMPA CERAMICAS E AZULEJOS LAVANDA 5L 1:15 START"; G=0;    H=$true  },
    @{ Row=6;  A=9; B="2025-04-28"; C=250; D="MAP SERVICOS DE CONSERVACAO - EIRELI";               E="000779"; F="PEDRA SANITARIA NAFT PLUS FLORAL 25G";                    G=116;  H=$false },
    @{ Row=7;  A=3; B="2025-04-30"; C=250; D="RH MULTI SERVICOS ADMINISTRATIVOS S.A";              E="000041"; F="LUVAS DESCARTAVEIS C/ 100 UND";                            G=1055; H=$false },
    @{ Row=8;  A=4; B="2025-05-07"; C=30;  D="V V REFEICOES LTDA";                                 E="000877"; F="CABO DE ALUMINIO NOBRE 140 CM COM PONTEIRA";               G=185;  H=$true  },
    @{ Row=9;  A=0; B="2025-05-08"; C=60;  D="ASSOCIACAO DOS EMPREENDEDORES DO SMVN";              E="000084"; F="SACO DE LIXO 200L PRETO 0.10 REFORCADO - PCT C/100 UND";  G=33;   H=$false },
    @{ Row=10; A=6; B="2025-05-08"; C=50;  D="ASSOCIACAO DOS EMPREENDEDORES DO SMVN";              E="000079"; F="SACO DE LIXO 50L REFORCADO - PCT C/100 UND";              G=-3;   H=$false },
    @{ Row=11; A=8; B="2025-05-08"; C=30;  D="CONDOMINIO SOBERANE RESIDENCE, CORPORATE E MALL";    E="000890"; F="AROMATIZANTE LIMPADOR PERF CONC COALA ALGODAO 120ML";     G=25;   H=$false }
)

function Set-TextValue($cell, $text) {
    # Temporarily force text format so Excel doesn't coerce a date-looking
    # or numeric-looking literal, then drop back to the default ("Normal")
    # cell style so no stray number-format is left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Write column by column (matching how the source workbook itself lays out
# its shared-string table) so new unique strings are introduced in the same
# relative order as the target file.
foreach ($r in $data) { $ws.Cells.Item($r.Row, 1).Value = $r.A }
foreach ($r in $data) { Set-TextValue $ws.Cells.Item($r.Row, 2) $r.B }
foreach ($r in $data) { $ws.Cells.Item($r.Row, 3).Value = $r.C }
foreach ($r in $data) { $ws.Cells.Item($r.Row, 4).Value = $r.D }
foreach ($r in $data) { Set-TextValue $ws.Cells.Item($r.Row, 5) $r.E }
foreach ($r in $data) { $ws.Cells.Item($r.Row, 6).Value = $r.F }
foreach ($r in $data) { $ws.Cells.Item($r.Row, 7).Value = $r.G }
foreach ($r in $data) { $ws.Cells.Item($r.Row, 8).Value = $r.H }

# Row 11 is brand new: give column A the same style as the rest of column A
# (s="1") by copying it from the cell above, then overwrite with its value.
$ws.Cells.Item(10, 1).Copy($ws.Cells.Item(11, 1))
$ws.Cells.Item(11, 1).Value = 8
